$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D+E updates ---
# Row 2: D 42.842.80 -> 42.840.74, E +0.32% -> +0.44%
$ws.Range("D2").Value = "'42.840.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "
# Row 3: D 2.570.27 -> 2.567.06, E +1.63% -> +1.50%
$ws.Range("D3").Value = "'2.567.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "
# Row 5: D 312.56 -> 311.80, E -0.95% -> -1.10%
$ws.Range("D5").Value = "'311.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
# Row 6: D 98.27 -> 98.49, E +2.88% -> +3.17%
$ws.Range("D6").Value = "'98.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.17%  "
# Row 9: D 0.533 -> 0.532, E +0.32% -> +0.36%
$ws.Range("D9").Value = "'0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
# Row 10: D 35.72 -> 35.64, E -0.13% -> -0.03%
$ws.Range("D10").Value = "'35.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
# Row 11: D 0.0809 -> 0.0808, E +0.60% -> +0.65%
$ws.Range("D11").Value = "'0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
# Row 12: D 7.44 -> 7.43, E -1.05% -> -0.98%
$ws.Range("D12").Value = "'7.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.98%  "
# Row 13: D 2.956.86 -> 2.958.15, E +1.29% -> +1.46%
$ws.Range("D13").Value = "'2.958.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "
# Row 15: D 15.83 -> 15.93, E +4.82% -> +6.01%
$ws.Range("D15").Value = "'15.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.01%  "
# Row 16: D 2.618.15 -> 2.513.39, E +2.59% -> +0.23%
$ws.Range("D16").Value = "'2.513.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
# Row 17: D 0.842 -> 0.840, E -0.48% -> -0.57%
$ws.Range("D17").Value = "'0.840"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.57%  "
# Row 18: D 42.858.47 -> 42.863.89, E +0.09% -> +0.27%
$ws.Range("D18").Value = "'42.863.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "
# Row 19: D 6.73 -> 6.71, E -1.70% -> -1.28%
$ws.Range("D19").Value = "'6.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.28%  "
# Row 22: D 69.74 -> 69.43, E +0.18% -> -0.09%
$ws.Range("D22").Value = "'69.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
# Row 23: D 248.64 -> 248.37, E -1.01% -> -0.80%
$ws.Range("D23").Value = "'248.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
# Row 26: D 27.05 -> 26.96, E +2.27% -> +2.25%
$ws.Range("D26").Value = "'26.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
# Row 28: D 2.41 -> 2.40, E -0.10% -> +0.15%
$ws.Range("D28").Value = "'2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
# Row 29: D 39.82 -> 39.89, E -1.28% -> -1.08%
$ws.Range("D29").Value = "'39.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
# Row 30: D 10.22 -> 10.18, E -1.56% -> -1.67%
$ws.Range("D30").Value = "'10.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "
# Row 31: D 159.57 -> 159.45, E +1.99% -> +2.11%
$ws.Range("D31").Value = "'159.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.11%  "
# Row 32: D 5.77 -> 5.76, E -1.97% -> -2.13%
$ws.Range("D32").Value = "'5.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.13%  "
# Row 36: D 3.30 -> 3.28, E -0.75% -> -0.55%
$ws.Range("D36").Value = "'3.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
# Row 37: D 18.65 -> 18.63, E -0.85% -> -0.56%
$ws.Range("D37").Value = "'18.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
# Row 38: D 2.57 -> 2.58, E +10.87% -> +13.29%
$ws.Range("D38").Value = "'2.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.29%  "
# Row 39: D 0.111 -> 0.112, E -0.12% -> +1.10%
$ws.Range("D39").Value = "'0.112"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
# Row 41: D 22.83 -> 22.86, E +2.55% -> +3.64%
$ws.Range("D41").Value = "'22.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.64%  "
# Row 42: D 4.11 -> 4.08, E +7.87% -> +7.48%
$ws.Range("D42").Value = "'4.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.48%  "
# Row 44: D 0.0302 -> 0.0301, E -0.27% -> -0.28%
$ws.Range("D44").Value = "'0.0301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "
# Row 48: D 2.809.09 -> 2.811.14, E +1.27% -> +1.46%
$ws.Range("D48").Value = "'2.811.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "
# Row 50: D 81.56 -> 81.23, E -3.44% -> -3.57%
$ws.Range("D50").Value = "'81.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.57%  "
# Row 51: D 74.23 -> 74.14, E -0.56% -> -0.93%
$ws.Range("D51").Value = "'74.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "

# --- E-only updates ---
# Row 4: E +0.08% -> +0.06%
$ws.Range("E4").Value = "  +0.06%  "
# Row 7: E +0.16% -> -0.06%
$ws.Range("E7").Value = "  -0.06%  "
# Row 8: E +0.00% -> -0.01%
$ws.Range("E8").Value = "  -0.01%  "
# Row 14: E -1.34% -> -1.38%
$ws.Range("E14").Value = "  -1.38%  "
# Row 24: E -0.23% -> +0.21%
$ws.Range("E24").Value = "  +0.21%  "
# Row 25: E -0.07% -> +0.48%
$ws.Range("E25").Value = "  +0.48%  "
# Row 27: E +0.01% -> -0.03%
$ws.Range("E27").Value = "  -0.03%  "
# Row 40: E -0.07% -> -0.22%
$ws.Range("E40").Value = "  -0.22%  "
# Row 43: E -0.11% -> -0.09%
$ws.Range("E43").Value = "  -0.09%  "
# Row 47: E -0.51% -> -0.37%
$ws.Range("E47").Value = "  -0.37%  "
# Row 49: E +2.71% -> +3.13%
$ws.Range("E49").Value = "  +3.13%  "

# --- Row 20/21 content swap (InternetComputer(DFINITY) <-> ShibaInu) ---
# Row 20 becomes ShibaInu
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0959"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
# Row 21 becomes InternetComputer(DFINITY)
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'12.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "

# --- Row 33/34/35 content shift (Hedera, ARBITRUM, WEMIXToken) ---
# Row 33 becomes ARBITRUM
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'2.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.89%  "
# Row 34 becomes WEMIXToken
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.10%  "
# Row 35 becomes Hedera
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0796"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.58%  "

# --- Row 45/46 content swap (NEARProtocol <-> Maker) ---
# Row 45 becomes Maker
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'1.992.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
# Row 46 becomes NEARProtocol
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "

